$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new value in the previously empty B5 cell
$ws.Range("B5").Value = "Redirecting to Plant"

# Fix the postcode-like values (remove the space)
$ws.Range("D2").Value = "NG156TF"
$ws.Range("D3").Value = "HD54LTF"

# Update the selected cell in the sheet view
$ws.Range("D3").Select()
